$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 239.5
$ws.Range("I19").Value = 254.2
$ws.Range("J19").Value = 232.81818
$ws.Range("K19").Value = 254.2
$ws.Range("L19").Value = 232.81818
$ws.Range("M19").Value = -79.19999999999999
$ws.Range("N19").Value = -582.81818

$ws.Range("H98").Value = 36896.91
$ws.Range("I98").Value = 934.46155
$ws.Range("J98").Value = 88842.664
$ws.Range("K98").Value = 934.46155
$ws.Range("L98").Value = 88842.664
$ws.Range("M98").Value = 563.53845
$ws.Range("N98").Value = -91838.664

$ws.Range("H117").Value = 48409.332
$ws.Range("J117").Value = 48409.332
$ws.Range("L117").Value = 48409.332
$ws.Range("N117").Value = -57587.332

$ws.Range("H122").Value = 36896.91
$ws.Range("I122").Value = 934.46155
$ws.Range("J122").Value = 88842.664
$ws.Range("K122").Value = 2803.38465
$ws.Range("L122").Value = 266527.992
$ws.Range("M122").Value = -353.38465
$ws.Range("N122").Value = -271427.992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 14980
$ws.Range("J11").Value = 14980
$ws.Range("L11").Value = 14980
$ws.Range("N11").Value = -15268

$ws.Range("H32").Value = 26847.865
$ws.Range("I32").Value = 30160.568
$ws.Range("K32").Value = 30160.568
$ws.Range("M32").Value = -29873.568

$ws.Range("H80").Value = 53317.168
$ws.Range("J80").Value = 53317.168
$ws.Range("L80").Value = 53317.168
$ws.Range("N80").Value = -55313.168

$ws.Range("H83").Value = 53317.168
$ws.Range("J83").Value = 53317.168
$ws.Range("L83").Value = 159951.504
$ws.Range("N83").Value = -169935.504

$ws.Range("H107").Value = 36971
$ws.Range("J107").Value = 36971
$ws.Range("L107").Value = 36971
$ws.Range("N107").Value = -44651

$ws.Range("H109").Value = 45092
$ws.Range("J109").Value = 45092
$ws.Range("L109").Value = 45092
$ws.Range("N109").Value = -47866

$ws.Range("H117").Value = 48412.2
$ws.Range("J117").Value = 48412.2
$ws.Range("L117").Value = 48412.2
$ws.Range("N117").Value = -57590.2

$ws.Range("H118").Value = 49626
$ws.Range("J118").Value = 49626
$ws.Range("L118").Value = 49626
$ws.Range("N118").Value = -52940

$ws.Range("H119").Value = 52684
$ws.Range("J119").Value = 52684
$ws.Range("L119").Value = 52684
$ws.Range("N119").Value = -62360

$ws.Range("H122").Value = 2324.32
$ws.Range("I122").Value = 2336.0952
$ws.Range("K122").Value = 7008.285600000001
$ws.Range("M122").Value = -4558.285600000001

$ws.Range("H131").Value = 48674.75
$ws.Range("J131").Value = 48674.75
$ws.Range("L131").Value = 48674.75
$ws.Range("N131").Value = -58754.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2030.3077
$ws.Range("I58").Value = 1730.05
$ws.Range("J58").Value = 3031.1667
$ws.Range("K58").Value = 1730.05
$ws.Range("L58").Value = 3031.1667
$ws.Range("M58").Value = -1527.05
$ws.Range("N58").Value = -3437.1667

$ws.Range("H99").Value = 1618.2727
$ws.Range("I99").Value = 1800.4
$ws.Range("J99").Value = 1466.5
$ws.Range("K99").Value = 1800.4
$ws.Range("L99").Value = 1466.5
$ws.Range("M99").Value = -302.4000000000001
$ws.Range("N99").Value = -4462.5

$ws.Range("H111").Value = 47650
$ws.Range("J111").Value = 47650
$ws.Range("L111").Value = 47650
$ws.Range("N111").Value = -55830

$ws.Range("H115").Value = 32069.5
$ws.Range("J115").Value = 32069.5
$ws.Range("L115").Value = 32069.5
$ws.Range("N115").Value = -34419.5

$ws.Range("H116").Value = 47891.5
$ws.Range("J116").Value = 47891.5
$ws.Range("L116").Value = 47891.5
$ws.Range("N116").Value = -57069.5

$ws.Range("H122").Value = 55517.047
$ws.Range("I122").Value = 67421.11
$ws.Range("J122").Value = 1948.75
$ws.Range("K122").Value = 202263.33
$ws.Range("L122").Value = 5846.25
$ws.Range("M122").Value = -199813.33
$ws.Range("N122").Value = -10746.25

$ws.Range("H126").Value = 1618.2727
$ws.Range("I126").Value = 1800.4
$ws.Range("J126").Value = 1466.5
$ws.Range("K126").Value = 5401.200000000001
$ws.Range("L126").Value = 4399.5
$ws.Range("M126").Value = -2931.200000000001
$ws.Range("N126").Value = -9339.5

$ws.Range("H136").Value = 2030.3077
$ws.Range("I136").Value = 1730.05
$ws.Range("J136").Value = 3031.1667
$ws.Range("K136").Value = 5190.15
$ws.Range("L136").Value = 9093.500100000001
$ws.Range("M136").Value = -2640.15
$ws.Range("N136").Value = -14193.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3200.3542
$ws.Range("I131").Value = 17203.166
$ws.Range("J131").Value = 1199.9524
$ws.Range("K131").Value = 51609.49800000001
$ws.Range("L131").Value = 3599.857199999999
$ws.Range("M131").Value = -46569.49800000001
$ws.Range("N131").Value = -13679.8572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 3110.8667
$ws.Range("I107").Value = 556.125
$ws.Range("J107").Value = 6030.5713
$ws.Range("K107").Value = 556.125
$ws.Range("L107").Value = 6030.5713
$ws.Range("M107").Value = 1363.875
$ws.Range("N107").Value = -9870.5713

$ws.Range("H122").Value = 1175.5883
$ws.Range("I122").Value = 1229.6154
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 3688.8462
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -1238.8462
$ws.Range("N122").Value = -7900

$ws.Range("H123").Value = 17487
$ws.Range("J123").Value = 17487
$ws.Range("L123").Value = 17487
$ws.Range("N123").Value = -22387

$ws.Range("H126").Value = 1660.619
$ws.Range("I126").Value = 2065.7273
$ws.Range("J126").Value = 1215
$ws.Range("K126").Value = 6197.1819
$ws.Range("L126").Value = 3645
$ws.Range("M126").Value = -3727.1819
$ws.Range("N126").Value = -8585

$ws.Range("H130").Value = 44860.89
$ws.Range("J130").Value = 44860.89
$ws.Range("L130").Value = 44860.89
$ws.Range("N130").Value = -54900.89

$ws.Range("H141").Value = 39336.637
$ws.Range("J141").Value = 39336.637
$ws.Range("L141").Value = 39336.637
$ws.Range("N141").Value = -49696.637

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3199.2856
$ws.Range("I40").Value = 2489.4443
$ws.Range("J40").Value = 4477
$ws.Range("K40").Value = 2489.4443
$ws.Range("L40").Value = 4477
$ws.Range("M40").Value = -2353.4443
$ws.Range("N40").Value = -4749

$ws.Range("H100").Value = 2324.5
$ws.Range("J100").Value = 2333.3333
$ws.Range("L100").Value = 2333.3333
$ws.Range("N100").Value = -3415.3333

$ws.Range("H110").Value = 45644
$ws.Range("J110").Value = 45644
$ws.Range("L110").Value = 45644
$ws.Range("N110").Value = -53824

$ws.Range("H122").Value = 2233.2173
$ws.Range("I122").Value = 2233.2173
$ws.Range("K122").Value = 6699.651899999999
$ws.Range("M122").Value = -4249.651899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 93338
$ws.Range("I9").Value = 80000
$ws.Range("K9").Value = 80000
$ws.Range("M9").Value = -79860

$ws.Range("H110").Value = 28513.6
$ws.Range("J110").Value = 28513.6
$ws.Range("L110").Value = 28513.6
$ws.Range("N110").Value = -36693.6

$ws.Range("H122").Value = 1311.7333
$ws.Range("I122").Value = 1288.2727
$ws.Range("K122").Value = 3864.8181
$ws.Range("M122").Value = -1414.8181

$ws.Range("H126").Value = 1803.1923
$ws.Range("I126").Value = 1518.9524
$ws.Range("J126").Value = 2997
$ws.Range("K126").Value = 4556.857199999999
$ws.Range("L126").Value = 8991
$ws.Range("M126").Value = -2086.857199999999
$ws.Range("N126").Value = -13931
